$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230493664741516
$ws.Range("B1").Value = 2.828868389129639
$ws.Range("C1").Value = 4.816621780395508
$ws.Range("D1").Value = 2.07628607749939
$ws.Range("E1").Value = 1.156054496765137
